$wb = $excel.ActiveWorkbook

# --- Sheet2 -> "reporting options": drop the old ~RFInput / RFSwitch example table ---
$wsRpt = $wb.Worksheets.Item("Sheet2")
$wsRpt.Range("B17:I21").Clear()
$wsRpt.Name = "reporting options"
$wsRpt.Range("B16").Select()

# --- TimePeriods: add a 5th "Def5" milestone-year column (F) ---
$wsTP = $wb.Worksheets.Item("TimePeriods")
$wsTP.Range("F27").Value = "Def5"
$wsTP.Range("F27").Interior.Color = 8454143
$wsTP.Range("F29").Value = 2005
$wsTP.Range("F30").Value = 2010
$wsTP.Range("F31").Value = 2020
$wsTP.Range("F32").Value = 2025
$wsTP.Range("F33").Value = 2026
$wsTP.Range("F34").Value = 2027
$wsTP.Range("F35").Value = 2028
$wsTP.Range("F36").Value = 2029
$wsTP.Range("F37").Value = 2030
$wsTP.Range("F38").Value = 2040
$wsTP.Range("F39").Value = 2050
$wsTP.Range("F29:F39").Interior.Color = 15495930
$wsTP.Range("F37").Select()

# --- Make "Region-Time Slices" the active sheet/tab (was Interpol_Extrapol_Defaults) ---
$wsRTS = $wb.Worksheets.Item("Region-Time Slices")
$wsRTS.Activate()
